# Auto-applied data refresh for Sheets (scheduled runner)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = "121.72727"
$ws.Range("I28").Value = "121.72727"
$ws.Range("K28").Value = "121.72727"
$ws.Range("M28").Value = "363.27273"

$ws.Range("H39").Value = "285.58334"
$ws.Range("I39").Value = "39.666668"
$ws.Range("J39").Value = "1023.3333"
$ws.Range("K39").Value = "119.000004"
$ws.Range("L39").Value = "3069.9999"
$ws.Range("M39").Value = "176.999996"
$ws.Range("N39").Value = "-3661.9999"

$ws.Range("H62").Value = "2628.5"
$ws.Range("I62").Value = "2610.625"
$ws.Range("J62").Value = "2700"
$ws.Range("K62").Value = "2610.625"
$ws.Range("L62").Value = "2700"
$ws.Range("M62").Value = "-1986.625"
$ws.Range("N62").Value = "-3948"

$ws.Range("H65").Value = "2628.5"
$ws.Range("I65").Value = "2610.625"
$ws.Range("J65").Value = "2700"
$ws.Range("K65").Value = "13053.125"
$ws.Range("L65").Value = "13500"
$ws.Range("M65").Value = "-9933.125"
$ws.Range("N65").Value = "-19740"

$ws.Range("H111").Value = "1183"
$ws.Range("I111").Value = "1400"
$ws.Range("J111").Value = "966"
$ws.Range("K111").Value = "4200"
$ws.Range("L111").Value = "2898"
$ws.Range("M111").Value = "-1133"
$ws.Range("N111").Value = "-9032"

$ws.Range("H129").Value = "112254.445"
$ws.Range("I129").Value = "572.5"
$ws.Range("J129").Value = "201600"
$ws.Range("K129").Value = "1717.5"
$ws.Range("L129").Value = "604800"
$ws.Range("M129").Value = "3282.5"
$ws.Range("N129").Value = "-614800"

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = "2385.0938"
$ws.Range("I32").Value = "2139"
$ws.Range("J32").Value = "4107.75"
$ws.Range("K32").Value = "2139"
$ws.Range("L32").Value = "4107.75"
$ws.Range("M32").Value = "-1852"
$ws.Range("N32").Value = "-4681.75"

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = "24518"
$ws.Range("J38").Value = "24518"
$ws.Range("L38").Value = "24518"
$ws.Range("N38").Value = "-25350"

$ws.Range("H94").Value = "1109.35"
$ws.Range("I94").Value = "790.4666999999999"
$ws.Range("J94").Value = "2066"
$ws.Range("K94").Value = "790.4666999999999"
$ws.Range("L94").Value = "2066"
$ws.Range("M94").Value = "-339.4666999999999"
$ws.Range("N94").Value = "-2968"

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H39").Value = "7320.4"
$ws.Range("I39").Value = "1534"
$ws.Range("K39").Value = "1534"
$ws.Range("M39").Value = "-1143"

$ws.Range("H49").Value = "7320.4"
$ws.Range("I49").Value = "1534"
$ws.Range("K49").Value = "1534"
$ws.Range("M49").Value = "-1352"

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = "886.125"
$ws.Range("I129").Value = "822.5"
$ws.Range("J129").Value = "949.75"
$ws.Range("K129").Value = "2467.5"
$ws.Range("L129").Value = "2849.25"
$ws.Range("M129").Value = "2532.5"
$ws.Range("N129").Value = "-12849.25"

$ws.Range("H131").Value = "6235108"
$ws.Range("J131").Value = "16016759"
$ws.Range("L131").Value = "48050277"
$ws.Range("N131").Value = "-48060357"

$ws.Range("H133").Value = "2899.375"
$ws.Range("I133").Value = "3032.5"
$ws.Range("J133").Value = "2500"
$ws.Range("K133").Value = "9097.5"
$ws.Range("L133").Value = "7500"
$ws.Range("M133").Value = "-4037.5"
$ws.Range("N133").Value = "-17620"

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = "8500"
$ws.Range("I80").Value = "4500"
$ws.Range("J80").Value = "10500"
$ws.Range("K80").Value = "4500"
$ws.Range("L80").Value = "10500"
$ws.Range("M80").Value = "-3502"
$ws.Range("N80").Value = "-12496"

$ws.Range("H83").Value = "8500"
$ws.Range("I83").Value = "4500"
$ws.Range("J83").Value = "10500"
$ws.Range("K83").Value = "22500"
$ws.Range("L83").Value = "52500"
$ws.Range("M83").Value = "-17508"
$ws.Range("N83").Value = "-62484"

$ws.Range("H97").Value = "633.96"
$ws.Range("I97").Value = "436.58334"
$ws.Range("J97").Value = "816.1539"
$ws.Range("K97").Value = "436.58334"
$ws.Range("L97").Value = "816.1539"
$ws.Range("M97").Value = "59.41665999999998"
$ws.Range("N97").Value = "-1808.1539"

$ws.Range("H107").Value = "636.25"
$ws.Range("I107").Value = "559.0625"
$ws.Range("J107").Value = "945"
$ws.Range("K107").Value = "559.0625"
$ws.Range("L107").Value = "945"
$ws.Range("M107").Value = "1360.9375"
$ws.Range("N107").Value = "-4785"

$ws.Range("H122").Value = "1962.0435"
$ws.Range("I122").Value = "1290.1111"
$ws.Range("J122").Value = "2394"
$ws.Range("K122").Value = "3870.3333"
$ws.Range("L122").Value = "7182"
$ws.Range("M122").Value = "-1420.3333"
$ws.Range("N122").Value = "-12082"

$ws.Range("H126").Value = "1965214.2"
$ws.Range("I126").Value = "5968.3335"
$ws.Range("J126").Value = "3033893.8"
$ws.Range("K126").Value = "17905.0005"
$ws.Range("L126").Value = "9101681.399999999"
$ws.Range("M126").Value = "-15435.0005"
$ws.Range("N126").Value = "-9106621.399999999"

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H38").Value = "0"
$ws.Range("I38").Value = "0"
$ws.Range("K38").Value = "0"
$ws.Range("M38").ClearContents()

$ws.Range("H132").Value = "14542.814"
$ws.Range("I132").Value = "22700.268"
$ws.Range("J132").Value = "4346"
$ws.Range("K132").Value = "68100.804"
$ws.Range("L132").Value = "13038"
$ws.Range("M132").Value = "-65570.804"
$ws.Range("N132").Value = "-18098"

$ws.Range("H136").Value = "5272.1113"
$ws.Range("I136").Value = "7689.8"
$ws.Range("J136").Value = "2250"
$ws.Range("K136").Value = "23069.4"
$ws.Range("L136").Value = "6750"
$ws.Range("M136").Value = "-20519.4"
$ws.Range("N136").Value = "-11850"

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = "55000"
$ws.Range("J16").Value = "55000"
$ws.Range("L16").Value = "55000"
$ws.Range("N16").Value = "-55584"

$ws.Range("H52").Value = "40047"
$ws.Range("J52").Value = "40047"
$ws.Range("L52").Value = "40047"
$ws.Range("N52").Value = "-40499"

$ws.Range("H107").Value = "548.5"
$ws.Range("I107").Value = "548.5"
$ws.Range("J107").Value = "0"
$ws.Range("K107").Value = "1645.5"
$ws.Range("L107").Value = "0"
$ws.Range("M107").Value = "274.5"
$ws.Range("N107").ClearContents()

$ws.Range("H132").Value = "2435.7144"
$ws.Range("I132").Value = "1614.5714"
$ws.Range("K132").Value = "4843.7142"
$ws.Range("M132").Value = "-2313.7142"
